# Restore C10 ("Integer min" for rule R20) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C10").Value = 1
